# Edit: Wed, Apr 01, 2020  9:05:43 AM
#
# 1) The table on slide 6 ("SOURCES OF FINANCE") gets its table style
#    switched from the deck's custom "Table_0" style to the built-in
#    "No Style, Table Grid" style ({0FD85157-5D88-4F22-A256-3A38BD79B148}).
# 2) The presentation's theme colour scheme is swapped from the
#    "Integral" palette over to the stock "Office Theme" palette
#    (dk1/lt1 stay black/white; dk2, lt2 and all six accents + the two
#    hyperlink colours change).

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 -------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{0FD85157-5D88-4F22-A256-3A38BD79B148}")

# --- 2) Swap the theme colour scheme to the "Office Theme" palette --
function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Item(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
